$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("emissions")
$ws2 = $wb.Worksheets.Item("removals")

# 1. Delete row 31 ("Steel, HRC (BF route)") on sheet1 -- shifts rows 32+ up by one
$ws1.Range("A31").EntireRow.Delete()

# 2. Add new note in J24 ("including disposal")
$ws1.Range("J24").Value = "including disposal"

# 3. Fix typo in sheet2 D4: "charcola (GLO)" -> "charcoal (GLO)"
$ws2.Range("D4").Value = "charcoal (GLO)"

# 4. Update view/selection state to match author's final position
$ws1.Activate()
$ws1.Range("J25").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 101
$ws2.Range("D4").Select()
